$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the port mapping for the "Webshell" web-server row (C13):
# previously "8082:82" (colliding with existing "8082:8082" entry) -> "8090:90"
$ws.Range("C13").Value = "8090:90"

# Move the active selection to F16 (was F7)
$ws.Range("F16").Select()
